$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "245.21"
Set-TextValue $ws.Range("G2") "22"

# Row 3
Set-TextValue $ws.Range("D3") "25.17"
Set-TextValue $ws.Range("G3") "22"

# Row 4
Set-TextValue $ws.Range("D4") "5.040"
Set-TextValue $ws.Range("G4") "22"

# Row 5
Set-TextValue $ws.Range("D5") "0.05610"
Set-TextValue $ws.Range("G5") "22"

# Row 6
Set-TextValue $ws.Range("D6") "6.554"
Set-TextValue $ws.Range("G6") "22"

# Row 7
Set-TextValue $ws.Range("D7") "3.020"
Set-TextValue $ws.Range("G7") "22"

# Row 8
Set-TextValue $ws.Range("D8") "0.8141"
Set-TextValue $ws.Range("G8") "22"

# Row 9
Set-TextValue $ws.Range("D9") "0.8425"
Set-TextValue $ws.Range("G9") "22"

# Row 10
Set-TextValue $ws.Range("D10") "0.1337"
Set-TextValue $ws.Range("G10") "22"

# Row 11
Set-TextValue $ws.Range("D11") "0.06959"
Set-TextValue $ws.Range("G11") "22"

# Row 12
Set-TextValue $ws.Range("D12") "0.02830"
Set-TextValue $ws.Range("G12") "22"

# Row 13
Set-TextValue $ws.Range("D13") "0.09411"
Set-TextValue $ws.Range("G13") "22"

# Row 14
Set-TextValue $ws.Range("D14") "0.001510"
Set-TextValue $ws.Range("G14") "22"

# Row 15
Set-TextValue $ws.Range("D15") "0.0005980"
Set-TextValue $ws.Range("G15") "22"

# Row 16
Set-TextValue $ws.Range("D16") "0.006245"
Set-TextValue $ws.Range("G16") "22"

# Row 17
Set-TextValue $ws.Range("G17") "22"

# Row 18
Set-TextValue $ws.Range("G18") "22"

# Row 19
Set-TextValue $ws.Range("D19") "0.3181"
Set-TextValue $ws.Range("G19") "22"

# Row 20
Set-TextValue $ws.Range("D20") "0.03264"
Set-TextValue $ws.Range("G20") "22"

# Row 21
Set-TextValue $ws.Range("D21") "0.1293"
Set-TextValue $ws.Range("G21") "22"

# Row 22
Set-TextValue $ws.Range("D22") "3.743"
Set-TextValue $ws.Range("G22") "22"

# Row 23
Set-TextValue $ws.Range("D23") "0.04676"
Set-TextValue $ws.Range("G23") "22"

# Row 24
Set-TextValue $ws.Range("D24") "0.1370"
Set-TextValue $ws.Range("G24") "22"

# Row 25
Set-TextValue $ws.Range("D25") "0.001242"
Set-TextValue $ws.Range("G25") "22"

# Row 26
Set-TextValue $ws.Range("D26") "0.004532"
Set-TextValue $ws.Range("G26") "22"

# Row 27
Set-TextValue $ws.Range("D27") "0.00009698"
Set-TextValue $ws.Range("G27") "22"

# Row 28
Set-TextValue $ws.Range("D28") "0.0001940"
Set-TextValue $ws.Range("G28") "22"

# Row 29
Set-TextValue $ws.Range("G29") "22"

# Row 30
Set-TextValue $ws.Range("G30") "22"

# Row 31
Set-TextValue $ws.Range("G31") "22"

# Row 32
Set-TextValue $ws.Range("G32") "22"

# Row 33
Set-TextValue $ws.Range("G33") "22"

# Row 34
Set-TextValue $ws.Range("G34") "22"

# Row 35
Set-TextValue $ws.Range("G35") "22"

# Row 36
Set-TextValue $ws.Range("G36") "22"

# Row 37
Set-TextValue $ws.Range("G37") "22"

# Row 38
Set-TextValue $ws.Range("G38") "22"

# Row 39
Set-TextValue $ws.Range("G39") "22"

# Row 40
Set-TextValue $ws.Range("D40") "0.03665"
Set-TextValue $ws.Range("G40") "22"

# Row 41
Set-TextValue $ws.Range("G41") "22"

# Row 42
Set-TextValue $ws.Range("B42") "CEJI"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.002717"
Set-TextValue $ws.Range("E42") "41CEJICEJI"
Set-TextValue $ws.Range("G42") "22"

# Row 43
Set-TextValue $ws.Range("B43") "KickToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.003382"
Set-TextValue $ws.Range("E43") "42KickTokenKICKWorstin24h"
Set-TextValue $ws.Range("G43") "22"

# Row 44
Set-TextValue $ws.Range("D44") "0.008192"
Set-TextValue $ws.Range("G44") "22"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005290"
Set-TextValue $ws.Range("G45") "22"

# Row 46
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("G46") "22"

# Row 47
Set-TextValue $ws.Range("D47") "0.2260"
Set-TextValue $ws.Range("E47") "46CoinbaseStockTokenCOIN"
Set-TextValue $ws.Range("G47") "22"

# Row 48
Set-TextValue $ws.Range("G48") "22"

# Row 49
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("G49") "22"

# Row 50
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("G50") "22"

# Row 51
Set-TextValue $ws.Range("G51") "22"

Write-Output "All updates applied"